$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 39 through 50 (old extra rows no longer tracked)
$ws.Range("A39:D50").ClearContents()

# Update row 37: World Cup Qualifiers (UEFA) / id -> SWPL / id
$ws.Range("C37").Value = "ceu82myq9rpq841ts3jl7uvis"
$ws.Range("B37").Value = "SWPL"

# Update row 38: WSL2 / id -> Eredivisie / id
$ws.Range("B38").Value = "Eredivisie"
$ws.Range("C38").Value = "aouykkl1rt7zo06sg0kbzkbh0"

# Restore the saved scroll/selection state
$ws.Range("C38").Select()
